$d = $word.ActiveDocument

$d.Content.Find.Execute("{%- if education %}", $true, $false, $false, $false, $false,
                         $true, 1, $false, "{% if education %}", 2)

$d.Content.Find.Execute("{%- if hobbies %}", $true, $false, $false, $false, $false,
                         $true, 1, $false, "{% if hobbies %}", 2)
